$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11 ("2021年") below the existing last row (row 10, "2020年").
# Column A holds the year label, formatted like the other year cells (copy
# style from A10 so it keeps the bold/centered/bordered look).
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").Value = "2021年"

# Numeric data for the 2021 row, one cell per industry column (B..AQ).
# E11 is intentionally left blank (no data reported for that industry).
$ws.Range("B11").Value = 4318.71
$ws.Range("C11").Value = 1164.35
$ws.Range("D11").Value = 156.98
$ws.Range("F11").Value = 1845.94
$ws.Range("G11").Value = 3435.62
$ws.Range("H11").Value = 296.06
$ws.Range("I11").Value = 1670.65
$ws.Range("J11").Value = 709.86
$ws.Range("K11").Value = 658
$ws.Range("L11").Value = 388.89
$ws.Range("M11").Value = 88.48999999999999
$ws.Range("N11").Value = 763.83
$ws.Range("O11").Value = 1923.37
$ws.Range("P11").Value = 108.07
$ws.Range("Q11").Value = 800.4
$ws.Range("R11").Value = 2915.6
$ws.Range("S11").Value = 89.18000000000001
$ws.Range("T11").Value = 3631.71
$ws.Range("U11").Value = 6.64
$ws.Range("V11").Value = 762.05
$ws.Range("W11").Value = 122.39
$ws.Range("X11").Value = 1085.26
$ws.Range("Y11").Value = 9163.42
$ws.Range("Z11").Value = 639.03
$ws.Range("AA11").Value = 844.1900000000001
$ws.Range("AB11").Value = 3.51
$ws.Range("AC11").Value = 72447.74000000001
$ws.Range("AD11").Value = 1926.35
$ws.Range("AE11").Value = 906.03
$ws.Range("AF11").Value = 7208.83
$ws.Range("AG11").Value = 4984.5
$ws.Range("AH11").Value = 800.33
$ws.Range("AI11").Value = 375.5
$ws.Range("AJ11").Value = 85.09999999999999
$ws.Range("AK11").Value = 4942.44
$ws.Range("AL11").Value = 985.54
$ws.Range("AM11").Value = 10149.46
$ws.Range("AN11").Value = 288.34
$ws.Range("AO11").Value = 707.0700000000001
$ws.Range("AP11").Value = 1288.59
$ws.Range("AQ11").Value = 207.18
